$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# The slide already contains shapes with ids 1,5,8,9,10,15,21,27,31,32.
# This sandboxed COM runtime assigns new shape ids as "smallest unused
# id >= 2", so to land the soon-to-be-created group on id 33 (matching
# the authored edit) we first burn through the ids that are still free
# below 33 by adding and immediately deleting throw-away shapes.
$existingIds = @(1, 5, 8, 9, 10, 15, 21, 27, 31, 32)
for ($id = 2; $id -lt 33; $id++) {
    if ($existingIds -notcontains $id) {
        $dummy = $s.Shapes.AddShape(1, 0, 0, 10, 10)
        $dummy.Delete()
    }
}

# Bring the ellipse (id 9) and the rectangle (id 10) to the front of the
# z-order so that, once grouped, the member shapes end up ordered exactly
# as in the target deck: 15, 21, 27, 8, 31, 32, 9, 10.
$shapeById9 = $null
$shapeById10 = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Id -eq 9) { $shapeById9 = $s.Shapes.Item($i) }
}
$shapeById9.ZOrder(0)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Id -eq 10) { $shapeById10 = $s.Shapes.Item($i) }
}
$shapeById10.ZOrder(0)

# Collect the eight shapes that make up the "edge size" diagram (the
# straight connectors, the rounded-rectangle outline, the two brackets,
# the circle and the square) and group them together.
$memberIndexes = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $id = $s.Shapes.Item($i).Id
    if ($id -eq 15 -or $id -eq 21 -or $id -eq 27 -or $id -eq 8 -or $id -eq 31 -or $id -eq 32 -or $id -eq 9 -or $id -eq 10) {
        $memberIndexes += $i
    }
}
$range = $s.Shapes.Range($memberIndexes)
$grp = $range.Group()
$grp.Name = "グループ化 32"
